$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 45 (old rows 45-72 shift down to 48-75)
$ws.Range("A45:A47").EntireRow.Insert()

# Row 45 (new block, Fecha = 2022-03-16)
$ws.Cells.Item(45,1).Value = 3
$ws.Cells.Item(45,2).Value = "Femacal de La Calera"
$ws.Cells.Item(45,3).Value = "Coquimbo"
$ws.Cells.Item(45,4).Value = 44636
$ws.Cells.Item(45,5).Value = 5
$ws.Cells.Item(45,6).Value = "Fruta"
$ws.Cells.Item(45,7).Value = 100104
$ws.Cells.Item(45,8).Value = "Frutos de pepita"
$ws.Cells.Item(45,9).Value = 100104003
$ws.Cells.Item(45,10).Value = "Membrillo"
$ws.Cells.Item(45,11).Value = "Champion"
$ws.Cells.Item(45,12).Value = "Especial"
$ws.Cells.Item(45,13).Value = 56
$ws.Cells.Item(45,14).Value = 15000
$ws.Cells.Item(45,15).Value = 15000
$ws.Cells.Item(45,16).Value = 15000
$ws.Cells.Item(45,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(45,18).Value = "Región de O'Higgins"
$ws.Cells.Item(45,19).Value = 833
$ws.Cells.Item(45,20).Value = 18
$ws.Cells.Item(45,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 46 (new block, Fecha = 2022-03-16)
$ws.Cells.Item(46,1).Value = 3
$ws.Cells.Item(46,2).Value = "Femacal de La Calera"
$ws.Cells.Item(46,3).Value = "Coquimbo"
$ws.Cells.Item(46,4).Value = 44636
$ws.Cells.Item(46,5).Value = 5
$ws.Cells.Item(46,6).Value = "Fruta"
$ws.Cells.Item(46,7).Value = 100104
$ws.Cells.Item(46,8).Value = "Frutos de pepita"
$ws.Cells.Item(46,9).Value = 100104003
$ws.Cells.Item(46,10).Value = "Membrillo"
$ws.Cells.Item(46,11).Value = "Champion"
$ws.Cells.Item(46,12).Value = "Primera"
$ws.Cells.Item(46,13).Value = 67
$ws.Cells.Item(46,14).Value = 13000
$ws.Cells.Item(46,15).Value = 13000
$ws.Cells.Item(46,16).Value = 13000
$ws.Cells.Item(46,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(46,18).Value = "Región de O'Higgins"
$ws.Cells.Item(46,19).Value = 722
$ws.Cells.Item(46,20).Value = 18
$ws.Cells.Item(46,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 47 (new block, Fecha = 2022-03-16)
$ws.Cells.Item(47,1).Value = 3
$ws.Cells.Item(47,2).Value = "Femacal de La Calera"
$ws.Cells.Item(47,3).Value = "Coquimbo"
$ws.Cells.Item(47,4).Value = 44636
$ws.Cells.Item(47,5).Value = 5
$ws.Cells.Item(47,6).Value = "Fruta"
$ws.Cells.Item(47,7).Value = 100104
$ws.Cells.Item(47,8).Value = "Frutos de pepita"
$ws.Cells.Item(47,9).Value = 100104003
$ws.Cells.Item(47,10).Value = "Membrillo"
$ws.Cells.Item(47,11).Value = "Champion"
$ws.Cells.Item(47,12).Value = "Segunda"
$ws.Cells.Item(47,13).Value = 60
$ws.Cells.Item(47,14).Value = 12000
$ws.Cells.Item(47,15).Value = 12000
$ws.Cells.Item(47,16).Value = 12000
$ws.Cells.Item(47,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(47,18).Value = "Región de O'Higgins"
$ws.Cells.Item(47,19).Value = 667
$ws.Cells.Item(47,20).Value = 18
$ws.Cells.Item(47,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append 3 new rows 73-75 at the end
# Row 73
$ws.Cells.Item(73,1).Value = 3
$ws.Cells.Item(73,2).Value = "Femacal de La Calera"
$ws.Cells.Item(73,3).Value = "Coquimbo"
$ws.Cells.Item(73,4).Value = 44628
$ws.Cells.Item(73,5).Value = 5
$ws.Cells.Item(73,6).Value = "Fruta"
$ws.Cells.Item(73,7).Value = 100104
$ws.Cells.Item(73,8).Value = "Frutos de pepita"
$ws.Cells.Item(73,9).Value = 100104003
$ws.Cells.Item(73,10).Value = "Membrillo"
$ws.Cells.Item(73,11).Value = "Champion"
$ws.Cells.Item(73,12).Value = "Especial"
$ws.Cells.Item(73,13).Value = 50
$ws.Cells.Item(73,14).Value = 15000
$ws.Cells.Item(73,15).Value = 15000
$ws.Cells.Item(73,16).Value = 15000
$ws.Cells.Item(73,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(73,18).Value = "Región de O'Higgins"
$ws.Cells.Item(73,19).Value = 833
$ws.Cells.Item(73,20).Value = 18
$ws.Cells.Item(73,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 74
$ws.Cells.Item(74,1).Value = 3
$ws.Cells.Item(74,2).Value = "Femacal de La Calera"
$ws.Cells.Item(74,3).Value = "Coquimbo"
$ws.Cells.Item(74,4).Value = 44628
$ws.Cells.Item(74,5).Value = 5
$ws.Cells.Item(74,6).Value = "Fruta"
$ws.Cells.Item(74,7).Value = 100104
$ws.Cells.Item(74,8).Value = "Frutos de pepita"
$ws.Cells.Item(74,9).Value = 100104003
$ws.Cells.Item(74,10).Value = "Membrillo"
$ws.Cells.Item(74,11).Value = "Champion"
$ws.Cells.Item(74,12).Value = "Primera"
$ws.Cells.Item(74,13).Value = 58
$ws.Cells.Item(74,14).Value = 13000
$ws.Cells.Item(74,15).Value = 13000
$ws.Cells.Item(74,16).Value = 13000
$ws.Cells.Item(74,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(74,18).Value = "Región de O'Higgins"
$ws.Cells.Item(74,19).Value = 722
$ws.Cells.Item(74,20).Value = 18
$ws.Cells.Item(74,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 75
$ws.Cells.Item(75,1).Value = 3
$ws.Cells.Item(75,2).Value = "Femacal de La Calera"
$ws.Cells.Item(75,3).Value = "Coquimbo"
$ws.Cells.Item(75,4).Value = 44628
$ws.Cells.Item(75,5).Value = 5
$ws.Cells.Item(75,6).Value = "Fruta"
$ws.Cells.Item(75,7).Value = 100104
$ws.Cells.Item(75,8).Value = "Frutos de pepita"
$ws.Cells.Item(75,9).Value = 100104003
$ws.Cells.Item(75,10).Value = "Membrillo"
$ws.Cells.Item(75,11).Value = "Champion"
$ws.Cells.Item(75,12).Value = "Segunda"
$ws.Cells.Item(75,13).Value = 56
$ws.Cells.Item(75,14).Value = 11000
$ws.Cells.Item(75,15).Value = 11000
$ws.Cells.Item(75,16).Value = 11000
$ws.Cells.Item(75,17).Value = "`$/caja 18 kilos empedrada"
$ws.Cells.Item(75,18).Value = "Región de O'Higgins"
$ws.Cells.Item(75,19).Value = 611
$ws.Cells.Item(75,20).Value = 18
$ws.Cells.Item(75,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
